# Extend the "group" columns on sheet1 from CL (last existing column) out to
# DG, replicating the per-row constant value that already fills columns C:CL.
# This mirrors the commit "Ajout de extract_column_from_all_sheets et d'un
# test." which widened the sample/test data table from A1:CL15 to A1:DG15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> value that already occupies C<row>:CL<row> and must be copied into
# the newly-appended CM<row>:DG<row> block. Row 10 is the blank test row.
$rowValues = @{
    2  = "group3"
    3  = "group1"
    4  = "group3"
    5  = "group3"
    6  = "group1"
    7  = "group1"
    8  = "group2"
    9  = "group2"
    10 = ""
    11 = "group2"
    12 = "group2"
    13 = "group2"
    14 = "group2"
    15 = "group1"
}

# Touch the bottom-right corner first so the sheet's used range/dimension
# grows to DG15 even though row 10's new cells end up blank.
$ws.Range("DG15").Value = $rowValues[15]

for ($r = 2; $r -le 15; $r++) {
    $value = $rowValues[$r]
    $targetRange = "CM" + $r + ":DG" + $r
    $sourceCell = "CL" + $r

    if ($value -eq "") {
        # Force the row's used range out to column DG, then blank it back out
        # so it matches the empty C10:CL10 cells already on this row.
        $ws.Range($targetRange).Value = "tmp"
        $ws.Range($targetRange).Value = ""
    } else {
        $ws.Range($targetRange).Value = $value
    }

    # Copy the formatting of the last existing cell in the row onto the new
    # block so the appended cells keep the same (default) style.
    $ws.Range($sourceCell).Copy()
    $ws.Range($targetRange).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
